$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.706.38"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.912.01"
$ws.Range("E3").Value = "  +0.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.20%  "

$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4915"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2957"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06733"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.16%  "

$ws.Range("D10").Value = "1.931.05"
$ws.Range("E10").Value = "  +1.82%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.47%  "

$ws.Range("E13").Value = "  +2.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6687"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "30.661.39"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.74%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007865"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

$ws.Range("D20").Value = "2.170.58"
$ws.Range("E20").Value = "  +1.37%  "

$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.280"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "194.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.243"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.61%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.613"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.940"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.478"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.339"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09086"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.053"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05222"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7352"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.107"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.731"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01811"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.713"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9176"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.060"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +28.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4427"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "106.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.897"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1379"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.552"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.065"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05863"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3982"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
